# TC09_Canine_Filter_Breed-AusShephd.xlsx - "Commiting to a stable branch"
#
# Populates the previously-empty query cell (B2) on the "startup" sheet with
# the Neo4j/Cypher query used to pull the filtered dataset, and brings the
# worksheet view/selection in line with that edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$query = 'MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN [''Australian Shepherd''] WITH DISTINCT c AS c, p, s, demo, diag RETURN coalesce(c.case_id,'''') AS `Case ID` , coalesce(s.clinical_study_designation,'''') AS `Study Code` , coalesce(s.clinical_study_type,'''') AS  `Study Type`, coalesce(demo.breed,'''') AS Breed , coalesce(diag.disease_term,'''') AS Diagnosis , coalesce(diag.stage_of_disease,'''') AS `Stage of Disease` ,  coalesce(demo.patient_age_at_enrollment,'''') AS Age , coalesce(demo.sex,'''') AS Sex , coalesce(demo.neutered_indicator,'''') AS  `Neutered Status`'

# B2 was blank (just carried the wrap-text style); give it the query text.
$ws.Range("B2").Value = $query

# With wrap-text on and ~930 chars of text, the row no longer fits at the
# default height, so it grows to accommodate the wrapped query.
$ws.Rows.Item(2).RowHeight = 188.5

# View scrolls left one column (topLeftCell C2 -> B2) now that B2 has content
# worth seeing, and the selection collapses from the whole column C down to
# just C2.
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("C2").Select()

Write-Output "Updated B2 with Cypher query and resized row 2"
